$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.308.62"
$ws.Range("E2").Value = "  -1.27%  "

$ws.Range("D3").Value = "2.772.49"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.548"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.90%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("E11").Value = "  +3.18%  "

$ws.Range("E12").Value = "  +3.28%  "

$ws.Range("E13").Value = "  -2.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "

$ws.Range("D15").Value = "3.206.40"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").Value = "2.768.36"
$ws.Range("E16").Value = "  -1.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.923"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("D18").Value = "51.256.22"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("E19").Value = "  +2.68%  "

$ws.Range("E20").Value = "  -2.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  -1.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.21%  "

$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.28%  "

$ws.Range("E28").Value = "  +13.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.39%  "

$ws.Range("E32").Value = "  +7.33%  "

$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0442"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0825"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("E39").Value = "  -2.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("E41").Value = "  -0.51%  "

$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("E45").Value = "  -1.99%  "

$ws.Range("D46").Value = "2.096.31"
$ws.Range("E46").Value = "  +1.34%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.33%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.18%  "

$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.902"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.36%  "
